# Generate Report for Archive
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for both tracked files/locales. Update every cell that
# shows that status (the per-locale status columns on the Overview sheet,
# and the Status column on each per-locale detail sheet) and shrink the
# status columns to fit the new (shorter) text, matching the column
# auto-fit Excel performs after a content change.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Columns resize to fit the new, shorter status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
